# Auto-update draw results: append the 2025-10-28 Pick 4 draw as a new row.
# Leading apostrophes force the date-like / digit-only values to be stored
# as literal text (matching the existing rows, which are all text cells),
# instead of being auto-converted to a date serial / number by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 42

$ws.Range("A$row").Value = "'2025-10-28"
$ws.Range("B$row").Value = "Pick 4"
$ws.Range("C$row").Value = "'251028"
$ws.Range("D$row").Value = "6-8-7-7"
$ws.Range("E$row").Value = "'2025-10-28T21:40:18.321+04:00"
